$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

$ws.Range("B10").Value = 'Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental desde a obtenção das amostras in situ até o preparo preliminar do material a ser analisado. Executar procedimentos de análises qualitativas bem como interpretar, avaliar e criticar os resultados obtidos. Objetivos Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Compreender e aplicar os procedimentos mais comuns de amostragem, coleta e preparação de amostras bem como os erros a não cometer no preparo das amostras de interesse ambiental. Empregar tratamentos preliminares no preparo das amostras: limpeza, secagem, moagem e peneiramento. Compreender as bases teóricas da química analítica qualitativa de interesse ambiental.'
$ws.Range("C10").Value = 'Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental desde a obtenção das amostras in situ até o preparo preliminar do material a ser analisado. Executar procedimentos de análises qualitativas bem como interpretar, avaliar e criticar os resultados obtidos. Objetivos Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Compreender e aplicar os procedimentos mais comuns de amostragem, coleta e preparação de amostras bem como os erros a não cometer no preparo das amostras de interesse ambiental. Empregar tratamentos preliminares no preparo das amostras: limpeza, secagem, moagem e peneiramento. Compreender as bases teóricas da química analítica qualitativa de interesse ambiental.'

$ws.Range("B11").Value = 'Train professionals at a higher level with the ability to know the sequence of chemical analysis procedures of environmental interest from obtaining samples in situ to the preliminary preparation of the material to be analyzed. Perform qualitative analysis procedures as well as interpret, evaluate and criticize the results obtained. Specific Objectives: To train professionals at a higher level with the ability to know the steps of the analytical sequence. Understand and apply the most common procedures for sampling, collecting and preparing samples, as well as errors not to make in preparing samples of environmental interest. Use preliminary treatments in sample preparation: cleaning, drying, grinding and sieving. Understand the theoretical bases of qualitative analytical chemistry of environmental interest'
$ws.Range("C11").Value = 'Train professionals at a higher level with the ability to know the sequence of chemical analysis procedures of environmental interest from obtaining samples in situ to the preliminary preparation of the material to be analyzed. Perform qualitative analysis procedures as well as interpret, evaluate and criticize the results obtained. Specific Objectives: To train professionals at a higher level with the ability to know the steps of the analytical sequence. Understand and apply the most common procedures for sampling, collecting and preparing samples, as well as errors not to make in preparing samples of environmental interest. Use preliminary treatments in sample preparation: cleaning, drying, grinding and sieving. Understand the theoretical bases of qualitative analytical chemistry of environmental interest'

$ws.Range("B13").Value = '7455355 - Robson da Silva Rocha'
$ws.Range("C13").Value = '7455355 - Robson da Silva Rocha'

$ws.Range("B14").Value = 'Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental. Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental.'
$ws.Range("C14").Value = 'Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental. Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental.'

$ws.Range("B15").Value = 'Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance.'
$ws.Range("C15").Value = 'Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance.'

$ws.Range("B16").Value = '- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Análise de sólidos, partículas, sedimentos.- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-).- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes'
$ws.Range("C16").Value = '- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Análise de sólidos, partículas, sedimentos.- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-).- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes'

$ws.Range("B17").Value = '- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples'
$ws.Range("C17").Value = '- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples'

$ws.Range("B19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'
$ws.Range("C19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'

$ws.Range("B20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'
$ws.Range("C20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'

$ws.Range("B21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado'
$ws.Range("C21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado'
